$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric-looking card number that must stay text.
# Format the cell as Text first so Excel doesn't coerce it to a number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 23.11.2023"

# --- Row 6 ---
$ws.Range("B6").Value = "25.11."
$ws.Range("C6").Value = "26.11."
$ws.Range("D6").Value = "KARTENZ./25.11 EDEKA RO"
$ws.Range("E6").Value = "149,98-"

# --- Row 7 ---
$ws.Range("B7").Value = "27.11."
$ws.Range("C7").Value = "28.11."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU WWVXWO"
$ws.Range("E7").Value = "90,85-"

# --- Row 8 ---
$ws.Range("B8").Value = "29.11."
$ws.Range("C8").Value = "30.11."
$ws.Range("D8").Value = "PAYPAL XKKPAE"
$ws.Range("E8").Value = "86,77-"

# --- Row 9 ---
$ws.Range("B9").Value = "01.12."
$ws.Range("C9").Value = "02.12."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-30840954"
$ws.Range("E9").Value = "53,56-"

# --- Row 10 ---
$ws.Range("B10").Value = "05.12."
$ws.Range("C10").Value = "06.12."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 61990789"
$ws.Range("E10").Value = "87,46-"

# --- Row 11: was an empty filler row, now gets a new transaction ---
$ws.Range("B11").Value = "08.12."
$ws.Range("C11").Value = "09.12."
$ws.Range("D11").Value = "RECHNUNG VODAFONE GMBH 55184783"
$ws.Range("E11").Value = "41,84-"
# E11 previously used the centered filler style (s=12); the populated
# amount cells in this table (E6:E10, E12) use the right-aligned style.
# Copy that formatting onto E11 so it matches the rest of the column.
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 12: closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 11.12.2023"
$ws.Range("E12").Value = "510,46-"

# --- Row 13: next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 16.12.2023"
